$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet / update its title to reflect the new "through" date
$ws.Name = "Through 2021-12-22"

# Row 9 (July) - 2021 columns (T/U/V)
$ws.Range("T9").Value = 9
$ws.Range("U9").Value = 140
$ws.Range("V9").Value = 0.0604

# Row 14 (December, through 12-21 -> 12-22)
$ws.Range("A14").Value = "December (through 12-22)"

$ws.Range("B14").Value = 4
$ws.Range("C14").Value = 26
$ws.Range("D14").Value = 0.1333

$ws.Range("F14").Value = 63
$ws.Range("G14").Value = 0.08699999999999999

$ws.Range("H14").Value = 10
$ws.Range("I14").Value = 78
$ws.Range("J14").Value = 0.1136

$ws.Range("L14").Value = 46
$ws.Range("M14").Value = 0.08

$ws.Range("O14").Value = 40
$ws.Range("P14").Value = 0.09089999999999999

$ws.Range("Q14").Value = 6
$ws.Range("R14").Value = 98
$ws.Range("S14").Value = 0.0577

$ws.Range("U14").Value = 143
$ws.Range("V14").Value = 0.0138

# Row 15 (Total)
$ws.Range("B15").Value = 37
$ws.Range("C15").Value = 284
$ws.Range("D15").Value = 0.1153

$ws.Range("F15").Value = 567
$ws.Range("G15").Value = 0.1028

$ws.Range("H15").Value = 73
$ws.Range("I15").Value = 836
$ws.Range("J15").Value = 0.0803

$ws.Range("L15").Value = 654
$ws.Range("M15").Value = 0.1066

$ws.Range("O15").Value = 520
$ws.Range("P15").Value = 0.1003

$ws.Range("Q15").Value = 70
$ws.Range("R15").Value = 1298
$ws.Range("S15").Value = 0.0512

$ws.Range("T15").Value = 102
$ws.Range("U15").Value = 1686
$ws.Range("V15").Value = 0.057
